# Apply the updated crypto price/volume snapshot values to Sheet1.
# Numeric-looking Price values are written with a leading apostrophe so
# Excel keeps them as literal text (matching the original inlineStr cells)
# instead of silently reformatting them as numbers (e.g. "1.000" -> 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.535.62"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.899.89"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'239.22"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.4893"
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("D8").Value = "'0.2916"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").Value = "'0.06663"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").Value = "1.890.68"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("D11").Value = "'16.88"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "'0.07318"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").Value = "'5.206"
$ws.Range("E13").Value = "  +3.73%  "
$ws.Range("D14").Value = "'88.63"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "'0.6636"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "30.501.22"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.000007841"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "'13.38"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'5.445"
$ws.Range("E20").Value = "  +14.98%  "
$ws.Range("D21").Value = "2.148.46"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'195.97"
$ws.Range("E23").Value = "  -7.83%  "
$ws.Range("D24").Value = "'6.151"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "'9.444"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").Value = "'162.12"
$ws.Range("E26").Value = "  +3.56%  "
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("D28").Value = "'1.929"
$ws.Range("E28").Value = "  +5.41%  "
$ws.Range("D29").Value = "'1.453"
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("D30").Value = "'4.318"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("D31").Value = "'0.09171"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").Value = "'4.142"
$ws.Range("E32").Value = "  +5.09%  "
$ws.Range("D33").Value = "'0.05209"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("D34").Value = "'0.7314"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").Value = "'1.107"
$ws.Range("E35").Value = "  +2.45%  "
$ws.Range("D36").Value = "'2.730"
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("D37").Value = "'0.01827"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("D38").Value = "'2.678"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").Value = "'0.9252"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").Value = "'2.049"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "'0.4377"
$ws.Range("E41").Value = "  -1.86%  "
$ws.Range("D42").Value = "'106.77"
$ws.Range("E42").Value = "  +1.70%  "
$ws.Range("D43").Value = "'5.891"
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("D44").Value = "'0.9947"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.1361"
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'67.69"
$ws.Range("E46").Value = "  +3.47%  "
$ws.Range("D47").Value = "'7.509"
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("D48").Value = "'8.977"
$ws.Range("E48").Value = "  +4.03%  "
$ws.Range("D49").Value = "'34.43"
$ws.Range("E49").Value = "  +3.35%  "
$ws.Range("D50").Value = "'0.05829"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("D51").Value = "'0.3911"
$ws.Range("E51").Value = "  -4.15%  "
